# Template "Yeu cau mua hang" (purchase request import template):
# remove the order-line "Don vi mua" (purchase unit) column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column K currently holds the "Chi tiet don hang / Don vi mua (*)"
# header (and the "Don vi" value on the sample data row). Deleting the
# whole column removes it and shifts every later column one place left.
$ws.Range("K1").EntireColumn.Delete()

# Park the selection on the cell that now occupies the freed position.
$ws.Range("K1").Select()
